$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 291 (pushes current rows 291:312 down to 293:314)
$ws.Rows("291:292").Insert()

# New weekly snapshot (date 2022-06-02 = serial 44714), same price bucket as the
# week it displaced (old row 291/292, now shifted to 293/294).
$ws.Cells.Item(291, 1).Value = 1
$ws.Cells.Item(291, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(291, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(291, 4).Value = 44714
$ws.Cells.Item(291, 5).Value = 15
$ws.Cells.Item(291, 6).Value = 100114014
$ws.Cells.Item(291, 7).Value = 'Betarraga'
$ws.Cells.Item(291, 8).Value = 'Sin especificar'
$ws.Cells.Item(291, 9).Value = 'Primera'
$ws.Cells.Item(291, 10).Value = 1000
$ws.Cells.Item(291, 11).Value = 450
$ws.Cells.Item(291, 12).Value = 500
$ws.Cells.Item(291, 13).Value = 475
$ws.Cells.Item(291, 14).Value = '$/paquete 4 unidades'
$ws.Cells.Item(291, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(291, 16).Value = 119
$ws.Cells.Item(291, 17).Value = 4
$ws.Cells.Item(291, 18).Value = 'Hortaliza'

$ws.Cells.Item(292, 1).Value = 1
$ws.Cells.Item(292, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(292, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(292, 4).Value = 44714
$ws.Cells.Item(292, 5).Value = 15
$ws.Cells.Item(292, 6).Value = 100114014
$ws.Cells.Item(292, 7).Value = 'Betarraga'
$ws.Cells.Item(292, 8).Value = 'Sin especificar'
$ws.Cells.Item(292, 9).Value = 'Segunda'
$ws.Cells.Item(292, 10).Value = 800
$ws.Cells.Item(292, 11).Value = 450
$ws.Cells.Item(292, 12).Value = 500
$ws.Cells.Item(292, 13).Value = 475
$ws.Cells.Item(292, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(292, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(292, 16).Value = 95
$ws.Cells.Item(292, 17).Value = 5
$ws.Cells.Item(292, 18).Value = 'Hortaliza'

Write-Host "Done"
